$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.590.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.57%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''3.147.29'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +2.92%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = '''  -0.10%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = '''564.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +2.66%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = '''144.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +3.77%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +0.12%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = '''3.141.44'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +2.88%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''0.497'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +2.03%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = '''6.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +5.31%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = '''  +1.56%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = '''0.467'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +2.29%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = '''36.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +3.29%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = '''  +1.80%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''3.653.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +2.83%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = '''64.614.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.42%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = '''  +1.30%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''3.146.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +3.07%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = '''516.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +6.09%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = '''  +4.04%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''14.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +3.29%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''0.718'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +5.00%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = '''7.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +4.23%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = '''12.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +3.76%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''79.14'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +1.42%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = '''0.997'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.29%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = '''8.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +15.83%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = '''  +5.12%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = '''2.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +3.82%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = '''  +0.19%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = '''26.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +3.15%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = '''  -0.56%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = '''  +2.16%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = '''549.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -6.16%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = '''  -0.09%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = '''  +3.50%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = '''  +4.33%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = '''0.0435'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +6.65%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = '''0.0826'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +4.29%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = '''3.158.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +7.95%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = '''0.122'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +3.32%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = '''dogwifhat'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = '''https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = '''2.76'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -2.62%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = '''Cosmos'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = '''8.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +1.42%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = '''0.266'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +9.81%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = '''  +7.17%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = '''1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.03%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = '''25.45'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +3.04%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = '''120.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +1.89%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = '''Stellar'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = '''0.109'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +0.41%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = '''PEPE'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = '''https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = '''0.0₃0519'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -1.55%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = '''  +2.77%  '
$ws.Range("E51").Style = "Normal"

